# Generate Report for Archive
#
# 1. The localization status for the zh-cn / de-de targets moved from
#    "Ready for handoff" to "In Translation" -- update it everywhere it
#    appears (Overview sheet status columns + each per-locale sheet's
#    Status column) so the text stays in sync across the workbook.
# 2. Because the new status text is shorter than the old one, the
#    status columns were re-sized (narrower) on the Overview sheet
#    (columns E & F) and on each per-locale sheet (column C).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text ---------------------------------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value     = "In Translation"
$dede.Range("C2").Value     = "In Translation"

# --- Narrow the status columns to fit the new, shorter text ---------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth     = 12.5
$dede.Columns.Item(3).ColumnWidth     = 12.5
